$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update hours worked on 2024-04-XX (row 51) from 1.5 to 2
$ws.Range("B51").Value = 2

$wb.Application.CalculateFull()
